$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. date style) from the row above, then overwrite values
$ws.Range("A13:H13").Copy($ws.Range("A14:H14"))

$ws.Range("A14").Value = 9686.69
$ws.Range("B14").Value = 9792.4500000000007
$ws.Range("C14").Value = 281.06
$ws.Range("D14").Value = 284.08999999999997
$ws.Range("E14").Value = $true
$ws.Range("F14").Value = 1.08
$ws.Range("G14").Value = 42620.766400462962
$ws.Range("H14").Value = $false
